# Correct SQL out filename variable
# Adds four new tracking rows (17-20) describing the DOCX_FILENAME_SPACES fix,
# fills in the previously-blank category cell on row 16, and tidies up a
# couple of redundant (visually no-op) cell styles along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 16 - fill in the missing "Animal dsRNA and -ssRNA (M) proposals"
#    category cell in column A, matching the black-font style already used
#    for this same label on row 6.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = $ws.Range("A6").Value()
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null

# Clear the redundant (visually identical) explicit-font style that had been
# applied to B16; it never changed B16's rendered appearance.
$ws.Range("B16").Style = "Normal"

# ---------------------------------------------------------------------------
# 2. Clean up a handful of other redundant styles elsewhere on the sheet
#    that likewise had zero visual effect (font 0 re-applied to itself).
# ---------------------------------------------------------------------------
$ws.Range("B9").Style = "Normal"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Style = "Normal"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Style = "Normal"

# D13/D15 used a wrap-text style that redundantly re-applied the default
# font; replace it with the equivalent wrap-text style already used by the
# rest of column D (e.g. D2) so the look stays identical.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 3. Append the four new rows documenting the DOCX_FILENAME_SPACES fix.
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "Plant virus (P) proposals"
$ws.Range("B17").Value = "2022.016P.A.Cilevirus_3ns[ ].docx"
$ws.Range("D17").Value = "delete space before period-docx"
$ws.Range("E17").Value = "DOCX_FILENAME_SPACES"

$ws.Range("A18").Value = "Animal dsRNA and -ssRNA (M) proposals"
$ws.Range("B18").Value = "2022.001M.Alpha[_]and[_]betanucleorhabdoviruses_6nsp.*"
$ws.Range("D18").Value = "replaced spaces with underscore"
$ws.Range("E18").Value = "DOCX_FILENAME_SPACES"

$ws.Range("A19").Value = "Bacterial viruses (B) proposals"
$ws.Range("B19").Value = "2022.003B.Abolish[_]Haartmanvirus.*"
$ws.Range("D19").Value = "replaced spaces with underscore"
$ws.Range("E19").Value = "DOCX_FILENAME_SPACES"

$ws.Range("A20").Value = "Plant virus (P) proposals"
$ws.Range("B20").Value = "2022.016P.A.v1.Cilevirus_3ns[ ].docx"
$ws.Range("D20").Value = "delete space before period-docx"
$ws.Range("E20").Value = "DOCX_FILENAME_SPACES"

# Column A on row 18 reuses the same highlighted black-font style as A6/A16.
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null

# Column D on the new rows wraps text like the rest of column D.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D17:D20").PasteSpecial(-4122) | Out-Null

# Column E gets a small Lucida Grande / black-text font (as pasted in from
# the source document) on all four new rows.
$ws.Range("E17").Font.Color = 0
$ws.Range("E17").Font.Name = "Lucida Grande"
$ws.Range("E17").Font.Size = 11
$ws.Range("E17").Copy() | Out-Null
$ws.Range("E18:E20").PasteSpecial(-4122) | Out-Null

# Match the source row heights (17pt, single line).
$ws.Rows.Item(17).RowHeight = 17
$ws.Rows.Item(18).RowHeight = 17
$ws.Rows.Item(19).RowHeight = 17
$ws.Rows.Item(20).RowHeight = 17

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Leave the selection on the last entry, as in the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("A20").Select() | Out-Null
